# Aggiornamento fino a 28 luglio
# Appends new daily rows (302-328) to the Ravarino COVID report sheet,
# extending the data range from A1:D301 to A1:D328.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: date serial (col A), nuovi pos. (col B), somma mobile 7gg. (col C),
# somma mobile 7gg. per 100mila abitanti (col D)
$data = @(
    @(44376, 0, 1, 16.17337861879347),
    @(44377, 0, 1, 16.17337861879347),
    @(44378, 0, 1, 16.17337861879347),
    @(44379, 0, 1, 16.17337861879347),
    @(44380, 0, 1, 16.17337861879347),
    @(44381, 0, 1, 16.17337861879347),
    @(44382, 0, 0, 0),
    @(44383, 0, 0, 0),
    @(44384, 1, 1, 16.17337861879347),
    @(44385, 1, 2, 32.34675723758694),
    @(44386, 0, 2, 32.34675723758694),
    @(44387, 0, 2, 32.34675723758694),
    @(44388, 0, 2, 32.34675723758694),
    @(44389, 0, 2, 32.34675723758694),
    @(44390, 0, 2, 32.34675723758694),
    @(44391, 0, 1, 16.17337861879347),
    @(44392, 0, 0, 0),
    @(44393, 0, 0, 0),
    @(44394, 0, 0, 0),
    @(44395, 0, 0, 0),
    @(44396, 0, 0, 0),
    @(44397, 0, 0, 0),
    @(44398, 0, 0, 0),
    @(44399, 0, 0, 0),
    @(44400, 2, 2, 32.34675723758694),
    @(44401, 3, 5, 80.86689309396733),
    @(44402, 2, 7, 113.2136503315543)
)

$firstNewRow = 302
$lastNewRow = $firstNewRow + $data.Count - 1

# Carry the date column's formatting (bold font, border, centered alignment,
# custom date number format) down from the last existing row onto the newly
# added date cells, just like dragging the fill handle down would.
$ws.Range("A301").Copy()
$ws.Range("A$firstNewRow`:A$lastNewRow").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$row = $firstNewRow
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}
